$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-10 with the new Sending/Target cluster combinations (FAPs/sCs/ECs x FAPs/sCs/ECs)
# following Dr Hou advice: add ECs as a sending cluster and recompute the LR-pair stats.

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Sema3e"
$ws.Cells.Item(2, 3).Value = "Plxnd1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.03595566666666666
$ws.Cells.Item(2, 8).Value = 0.107867
$ws.Cells.Item(2, 9).Value = 0.01609691666901703
$ws.Cells.Item(2, 10).Value = 0.01609691666901703
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 73.19890333333333
$ws.Cells.Item(2, 14).Value = 219.59671
$ws.Cells.Item(2, 15).Value = 0.6596328743217019
$ws.Cells.Item(2, 16).Value = 0.6596328743217019
$ws.Cells.Item(2, 17).Value = 2.631915368618889
$ws.Cells.Item(2, 18).Value = 23.68723831757
$ws.Cells.Item(2, 19).Value = 0.01061805541010062
$ws.Cells.Item(2, 20).Value = 0.01061805541010062

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Sema3e"
$ws.Cells.Item(3, 3).Value = "Plxnd1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.03595566666666666
$ws.Cells.Item(3, 8).Value = 0.107867
$ws.Cells.Item(3, 9).Value = 0.01609691666901703
$ws.Cells.Item(3, 10).Value = 0.01609691666901703
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 9.661646333333335
$ws.Cells.Item(3, 14).Value = 28.984939
$ws.Cells.Item(3, 15).Value = 0.08706605224007774
$ws.Cells.Item(3, 16).Value = 0.08706605224007773
$ws.Cells.Item(3, 17).Value = 0.3473909350125556
$ws.Cells.Item(3, 18).Value = 3.126518415113
$ws.Cells.Item(3, 19).Value = 0.001401494987608814
$ws.Cells.Item(3, 20).Value = 0.001401494987608814

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Sema3e"
$ws.Cells.Item(4, 3).Value = "Plxnd1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.03595566666666666
$ws.Cells.Item(4, 8).Value = 0.107867
$ws.Cells.Item(4, 9).Value = 0.01609691666901703
$ws.Cells.Item(4, 10).Value = 0.01609691666901703
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 28.10860633333333
$ws.Cells.Item(4, 14).Value = 84.325819
$ws.Cells.Item(4, 15).Value = 0.2533010734382203
$ws.Cells.Item(4, 16).Value = 0.2533010734382204
$ws.Cells.Item(4, 17).Value = 1.010663679785889
$ws.Cells.Item(4, 18).Value = 9.095973118073
$ws.Cells.Item(4, 19).Value = 0.004077366271307594
$ws.Cells.Item(4, 20).Value = 0.004077366271307596

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Sema3e"
$ws.Cells.Item(5, 3).Value = "Plxnd1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.375637
$ws.Cells.Item(5, 8).Value = 4.126911
$ws.Cells.Item(5, 9).Value = 0.61585603073646
$ws.Cells.Item(5, 10).Value = 0.6158560307364601
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 73.19890333333333
$ws.Cells.Item(5, 14).Value = 219.59671
$ws.Cells.Item(5, 15).Value = 0.6596328743217019
$ws.Cells.Item(5, 16).Value = 0.6596328743217019
$ws.Cells.Item(5, 17).Value = 100.6951197847567
$ws.Cells.Item(5, 18).Value = 906.2560780628099
$ws.Cells.Item(5, 19).Value = 0.4062388837230455
$ws.Cells.Item(5, 20).Value = 0.4062388837230456

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Sema3e"
$ws.Cells.Item(6, 3).Value = "Plxnd1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.375637
$ws.Cells.Item(6, 8).Value = 4.126911
$ws.Cells.Item(6, 9).Value = 0.61585603073646
$ws.Cells.Item(6, 10).Value = 0.6158560307364601
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.661646333333335
$ws.Cells.Item(6, 14).Value = 28.984939
$ws.Cells.Item(6, 15).Value = 0.08706605224007774
$ws.Cells.Item(6, 16).Value = 0.08706605224007773
$ws.Cells.Item(6, 17).Value = 13.29091817704767
$ws.Cells.Item(6, 18).Value = 119.618263593429
$ws.Cells.Item(6, 19).Value = 0.05362015334446755
$ws.Cells.Item(6, 20).Value = 0.05362015334446755

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Sema3e"
$ws.Cells.Item(7, 3).Value = "Plxnd1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.375637
$ws.Cells.Item(7, 8).Value = 4.126911
$ws.Cells.Item(7, 9).Value = 0.61585603073646
$ws.Cells.Item(7, 10).Value = 0.6158560307364601
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 28.10860633333333
$ws.Cells.Item(7, 14).Value = 84.325819
$ws.Cells.Item(7, 15).Value = 0.2533010734382203
$ws.Cells.Item(7, 16).Value = 0.2533010734382204
$ws.Cells.Item(7, 17).Value = 38.66723889056766
$ws.Cells.Item(7, 18).Value = 348.0051500151089
$ws.Cells.Item(7, 19).Value = 0.1559969936689469
$ws.Cells.Item(7, 20).Value = 0.155996993668947

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Sema3e"
$ws.Cells.Item(8, 3).Value = "Plxnd1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8221063333333333
$ws.Cells.Item(8, 8).Value = 2.466319
$ws.Cells.Item(8, 9).Value = 0.3680470525945229
$ws.Cells.Item(8, 10).Value = 0.368047052594523
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 73.19890333333333
$ws.Cells.Item(8, 14).Value = 219.59671
$ws.Cells.Item(8, 15).Value = 0.6596328743217019
$ws.Cells.Item(8, 16).Value = 0.6596328743217019
$ws.Cells.Item(8, 17).Value = 60.17728202338778
$ws.Cells.Item(8, 18).Value = 541.59553821049
$ws.Cells.Item(8, 19).Value = 0.2427759351885558
$ws.Cells.Item(8, 20).Value = 0.2427759351885558

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Sema3e"
$ws.Cells.Item(9, 3).Value = "Plxnd1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8221063333333333
$ws.Cells.Item(9, 8).Value = 2.466319
$ws.Cells.Item(9, 9).Value = 0.3680470525945229
$ws.Cells.Item(9, 10).Value = 0.368047052594523
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 9.661646333333335
$ws.Cells.Item(9, 14).Value = 28.984939
$ws.Cells.Item(9, 15).Value = 0.08706605224007774
$ws.Cells.Item(9, 16).Value = 0.08706605224007773
$ws.Cells.Item(9, 17).Value = 7.942900641060112
$ws.Cells.Item(9, 18).Value = 71.48610576954101
$ws.Cells.Item(9, 19).Value = 0.03204440390800137
$ws.Cells.Item(9, 20).Value = 0.03204440390800137

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Sema3e"
$ws.Cells.Item(10, 3).Value = "Plxnd1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8221063333333333
$ws.Cells.Item(10, 8).Value = 2.466319
$ws.Cells.Item(10, 9).Value = 0.3680470525945229
$ws.Cells.Item(10, 10).Value = 0.368047052594523
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 28.10860633333333
$ws.Cells.Item(10, 14).Value = 84.325819
$ws.Cells.Item(10, 15).Value = 0.2533010734382203
$ws.Cells.Item(10, 16).Value = 0.2533010734382204
$ws.Cells.Item(10, 17).Value = 23.10826328780677
$ws.Cells.Item(10, 18).Value = 207.974369590261
$ws.Cells.Item(10, 19).Value = 0.0932267134979658
$ws.Cells.Item(10, 20).Value = 0.09322671349796584
